$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.616.85"
$ws.Range("E2").Value = "  +0.74%  "
$ws.Range("D3").Value = "2.792.87"
$ws.Range("E3").Value = "  +1.56%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "352.37"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.78%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "111.36"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.32%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.555"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.78%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  +8.44%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.96"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.01%  "
$ws.Range("E11").Value = "  -1.23%  "
$ws.Range("E12").Value = "  +0.05%  "
$ws.Range("E13").Value = "  +1.85%  "
$ws.Range("E14").Value = "  +3.28%  "
$ws.Range("D15").Value = "3.233.12"
$ws.Range("E15").Value = "  +1.27%  "
$ws.Range("D16").Value = "2.808.75"
$ws.Range("E16").Value = "  +2.06%  "
$ws.Range("E17").Value = "  +2.26%  "
$ws.Range("D18").Value = "51.578.46"
$ws.Range("E18").Value = "  +0.80%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.57"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.75%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.18"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +6.18%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.50"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.04%  "
$ws.Range("D22").Value = "0.0₃0968"
$ws.Range("E22").Value = "  +1.20%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.14"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.06%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "267.03"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.03%  "
$ws.Range("E25").Value = "  +0.58%  "
$ws.Range("E26").Value = "  -0.07%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "26.03"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.37%  "
$ws.Range("E28").Value = "  -1.72%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "38.92"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +11.58%  "
$ws.Range("E30").Value = "  +2.41%  "
$ws.Range("E31").Value = "  +0.51%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "52.60"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.01%  "
$ws.Range("E33").Value = "  +0.54%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0453"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.13%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0888"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +7.01%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.55"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +7.94%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.01%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.76"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.44%  "
$ws.Range("E39").Value = "  +3.51%  "
$ws.Range("E40").Value = "  +0.92%  "
$ws.Range("E41").Value = "  +1.45%  "
$ws.Range("E42").Value = "  +0.47%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.23"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.26%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "120.41"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.29%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "21.66"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.14%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.43"
$ws.Range("D46").Style = "Normal"
$ws.Range("E47").Value = "  +5.42%  "
$ws.Range("D48").Value = "2.102.89"
$ws.Range("E48").Value = "  +1.36%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.956"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.05%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.47"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.23%  "
$ws.Range("E51").Value = "  +6.88%  "
